$wb = $excel.ActiveWorkbook

# Column F ("想去人数" / "want-to-go count") updates, applied identically
# to both the "展览" and "全部类型" sheets (their data is duplicated).
$updates = @{
    2  = 165
    3  = 425
    4  = 12338
    6  = 139
    9  = 158
    11 = 447
    15 = 41
    16 = 368
    17 = 3361
    20 = 16
    22 = 36
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
